# Update "想去人数" (interested-count) figures with freshly scraped numbers.
# Source: gh-pages output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3290   # CW国潮动漫游戏嘉年华
$wsExpo.Range("F5").Value = 6913   # 第二届漫画城市动漫展 -故事再次开始
$wsExpo.Range("F6").Value = 2211   # 环形宇宙动漫游戏嘉年华
$wsExpo.Range("F8").Value = 86     # 环形宇宙动漫游戏嘉年华内场票-钱文青
$wsExpo.Range("F11").Value = 0     # 第二届漫画城市动漫展内场-《琅声雅集》

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 20     # 全国地下偶像联合公演展-永乐大典Vol.01

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3290    # CW国潮动漫游戏嘉年华
$wsAll.Range("F3").Value = 20      # 全国地下偶像联合公演展-永乐大典Vol.01
$wsAll.Range("F6").Value = 6913    # 第二届漫画城市动漫展 -故事再次开始
$wsAll.Range("F7").Value = 2211    # 环形宇宙动漫游戏嘉年华
$wsAll.Range("F9").Value = 86      # 环形宇宙动漫游戏嘉年华内场票-钱文青
$wsAll.Range("F12").Value = 74     # 第二届漫画城市动漫展内场-《琅声雅集》
